$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing row for 2020-03-25 (regional data that was missing a link)
$row = 24
$ws.Cells.Item($row, 1).Value = 43915   # fecha
$ws.Cells.Item($row, 2).Value = 23      # dia
$ws.Cells.Item($row, 3).Value = 3       # Arica y Parinacota
$ws.Cells.Item($row, 4).Value = 5       # Tarapacá
$ws.Cells.Item($row, 5).Value = 20      # Antofagasta
$ws.Cells.Item($row, 6).Value = 1       # Atacama
$ws.Cells.Item($row, 7).Value = 13      # Coquimbo
$ws.Cells.Item($row, 8).Value = 44      # Valparaíso
$ws.Cells.Item($row, 9).Value = 746     # Metropolitana
$ws.Cells.Item($row, 10).Value = 14     # O'Higgins
$ws.Cells.Item($row, 11).Value = 31     # Maule
$ws.Cells.Item($row, 12).Value = 114    # Ñuble
$ws.Cells.Item($row, 13).Value = 109    # Biobío
$ws.Cells.Item($row, 14).Value = 111    # Araucanía
$ws.Cells.Item($row, 15).Value = 14     # Los Ríos
$ws.Cells.Item($row, 16).Value = 60     # Los Lagos
$ws.Cells.Item($row, 17).Value = 2      # Aysén
$ws.Cells.Item($row, 18).Value = 19     # Magallanes
$ws.Cells.Item($row, 19).Value = 1306   # total

# Match the fecha column date-number style used by the rest of column A
$ws.Cells.Item($row, 1).NumberFormat = "DD/MM/YY"

# Update the view state to match the diff (scroll + selection)
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("K57").Select()
